$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 350.34784
$ws.Range("I33").Value = 307.5238
$ws.Range("K33").Value = 307.5238
$ws.Range("M33").Value = -78.52379999999999
$ws.Range("H129").Value = 770.8570999999999
$ws.Range("J129").Value = 899.1111
$ws.Range("L129").Value = 2697.3333
$ws.Range("N129").Value = -12697.3333
$ws.Range("H137").Value = 20879.057
$ws.Range("I137").Value = 2311.5518
$ws.Range("J137").Value = 43314.793
$ws.Range("K137").Value = 6934.655400000001
$ws.Range("L137").Value = 129944.379
$ws.Range("M137").Value = -4384.655400000001
$ws.Range("N137").Value = -135044.379
$ws.Range("H141").Value = 2470.9092
$ws.Range("I141").Value = 1755.7894
$ws.Range("J141").Value = 7000
$ws.Range("K141").Value = 5267.3682
$ws.Range("L141").Value = 21000
$ws.Range("M141").Value = -87.36819999999989
$ws.Range("N141").Value = -31360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15409.627
$ws.Range("I32").Value = 15852.972
$ws.Range("K32").Value = 15852.972
$ws.Range("M32").Value = -15565.972
$ws.Range("H132").Value = 15553.583
$ws.Range("I132").Value = 1599.3529
$ws.Range("J132").Value = 28038.947
$ws.Range("K132").Value = 4798.0587
$ws.Range("L132").Value = 84116.841
$ws.Range("M132").Value = -2268.0587
$ws.Range("N132").Value = -89176.841

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 1000
$ws.Range("M107").Value = 920
$ws.Range("H134").Value = 27984.385
$ws.Range("I134").Value = 30205.305
$ws.Range("J134").Value = 1333.3334
$ws.Range("K134").Value = 90615.91500000001
$ws.Range("L134").Value = 4000.0002
$ws.Range("M134").Value = -88080.91500000001
$ws.Range("N134").Value = -9070.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11783.452
$ws.Range("I31").Value = 23266.555
$ws.Range("J31").Value = 3171.125
$ws.Range("K31").Value = 23266.555
$ws.Range("L31").Value = 3171.125
$ws.Range("M31").Value = -22971.555
$ws.Range("N31").Value = -3761.125
$ws.Range("H34").Value = 11783.452
$ws.Range("I34").Value = 23266.555
$ws.Range("J34").Value = 3171.125
$ws.Range("K34").Value = 23266.555
$ws.Range("L34").Value = 3171.125
$ws.Range("M34").Value = -23064.555
$ws.Range("N34").Value = -3575.125
$ws.Range("H132").Value = 23948.379
$ws.Range("I132").Value = 32083.85
$ws.Range("J132").Value = 5869.5557
$ws.Range("K132").Value = 96251.54999999999
$ws.Range("L132").Value = 17608.6671
$ws.Range("M132").Value = -93721.54999999999
$ws.Range("N132").Value = -22668.6671
$ws.Range("H134").Value = 3331.2104
$ws.Range("I134").Value = 605.4706
$ws.Range("J134").Value = 26500
$ws.Range("K134").Value = 1816.4118
$ws.Range("L134").Value = 79500
$ws.Range("M134").Value = 718.5882000000001
$ws.Range("N134").Value = -84570

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 5150.12
$ws.Range("I68").Value = 769.8
$ws.Range("J68").Value = 6245.2
$ws.Range("K68").Value = 2309.4
$ws.Range("L68").Value = 18735.6
$ws.Range("M68").Value = -1498.4
$ws.Range("N68").Value = -20357.6
$ws.Range("H71").Value = 5150.12
$ws.Range("I71").Value = 769.8
$ws.Range("J71").Value = 6245.2
$ws.Range("K71").Value = 6928.2
$ws.Range("L71").Value = 56206.8
$ws.Range("M71").Value = -2872.2
$ws.Range("N71").Value = -64318.8
$ws.Range("H131").Value = 806.09
$ws.Range("J131").Value = 815.5567
$ws.Range("L131").Value = 2446.6701
$ws.Range("N131").Value = -12526.6701
$ws.Range("H132").Value = 799.8570999999999
$ws.Range("I132").Value = 799.8333
$ws.Range("K132").Value = 7198.4997
$ws.Range("M132").Value = -4668.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10739.571
$ws.Range("I80").Value = 22841
$ws.Range("J80").Value = 4016.5557
$ws.Range("K80").Value = 22841
$ws.Range("L80").Value = 4016.5557
$ws.Range("M80").Value = -21843
$ws.Range("N80").Value = -6012.5557
$ws.Range("H83").Value = 10739.571
$ws.Range("I83").Value = 22841
$ws.Range("J83").Value = 4016.5557
$ws.Range("K83").Value = 114205
$ws.Range("L83").Value = 20082.7785
$ws.Range("M83").Value = -109213
$ws.Range("N83").Value = -30066.7785
$ws.Range("H113").Value = 4720.3335
$ws.Range("I113").Value = 4750
$ws.Range("J113").Value = 4661
$ws.Range("K113").Value = 4750
$ws.Range("L113").Value = 4661
$ws.Range("M113").Value = -2580
$ws.Range("N113").Value = -9001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3587.2144
$ws.Range("I68").Value = 2306.5
$ws.Range("J68").Value = 4099.5
$ws.Range("K68").Value = 2306.5
$ws.Range("L68").Value = 4099.5
$ws.Range("M68").Value = -1557.5
$ws.Range("N68").Value = -5597.5
$ws.Range("H71").Value = 3587.2144
$ws.Range("I71").Value = 2306.5
$ws.Range("J71").Value = 4099.5
$ws.Range("K71").Value = 11532.5
$ws.Range("L71").Value = 20497.5
$ws.Range("M71").Value = -7788.5
$ws.Range("N71").Value = -27985.5
$ws.Range("H110").Value = 2529724.8
$ws.Range("J110").Value = 2529724.8
$ws.Range("L110").Value = 2529724.8
$ws.Range("N110").Value = -2537904.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4425.75
$ws.Range("I62").Value = 3500
$ws.Range("J62").Value = 4734.3335
$ws.Range("K62").Value = 3500
$ws.Range("L62").Value = 4734.3335
$ws.Range("M62").Value = -2876
$ws.Range("N62").Value = -5982.3335
$ws.Range("H65").Value = 4425.75
$ws.Range("I65").Value = 3500
$ws.Range("J65").Value = 4734.3335
$ws.Range("K65").Value = 17500
$ws.Range("L65").Value = 23671.6675
$ws.Range("M65").Value = -14380
$ws.Range("N65").Value = -29911.6675
$ws.Range("H81").Value = 1792.375
$ws.Range("I81").Value = 1296.5714
$ws.Range("J81").Value = 5263
$ws.Range("K81").Value = 2593.1428
$ws.Range("L81").Value = 10526
$ws.Range("M81").Value = -1532.1428
$ws.Range("N81").Value = -12648
$ws.Range("H84").Value = 1792.375
$ws.Range("I84").Value = 1296.5714
$ws.Range("J84").Value = 5263
$ws.Range("K84").Value = 12965.714
$ws.Range("L84").Value = 52630
$ws.Range("M84").Value = -7661.714
$ws.Range("N84").Value = -63238
$ws.Range("H123").Value = 40429
$ws.Range("J123").Value = 40429
$ws.Range("L123").Value = 40429
$ws.Range("N123").Value = -50229
$ws.Range("H132").Value = 2112.0908
$ws.Range("I132").Value = 1906.9
$ws.Range("J132").Value = 2551.7856
$ws.Range("K132").Value = 5720.700000000001
$ws.Range("L132").Value = 7655.3568
$ws.Range("M132").Value = -3190.700000000001
$ws.Range("N132").Value = -12715.3568
$ws.Range("H136").Value = 1199.4667
$ws.Range("I136").Value = 699.36365
$ws.Range("J136").Value = 2574.75
$ws.Range("K136").Value = 2098.09095
$ws.Range("L136").Value = 7724.25
$ws.Range("M136").Value = 451.9090500000002
$ws.Range("N136").Value = -12824.25
